$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to remain text so numeric-looking strings
# like "67.208.83" or "0.999" are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "67.208.83"
$ws.Range("E2").Value = "  -1.66%  "
$ws.Range("D3").Value = "2.488.66"
$ws.Range("E3").Value = "  -2.00%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("D5").Value = "586.19"
$ws.Range("E5").Value = "  -1.46%  "
$ws.Range("D6").Value = "167.95"
$ws.Range("E6").Value = "  -5.57%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "0.518"
$ws.Range("E8").Value = "  -2.67%  "
$ws.Range("D9").Value = "2.489.10"
$ws.Range("E9").Value = "  -1.94%  "
$ws.Range("E10").Value = "  -3.96%  "
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").Value = "0.341"
$ws.Range("E12").Value = "  -2.25%  "
$ws.Range("E13").Value = "  -3.99%  "
$ws.Range("E14").Value = "  -3.75%  "
$ws.Range("E15").Value = "  -2.39%  "
$ws.Range("E16").Value = "  -3.56%  "
$ws.Range("D17").Value = "66.924.81"
$ws.Range("E17").Value = "  -1.48%  "
$ws.Range("D18").Value = "2.514.71"
$ws.Range("E18").Value = "  -1.50%  "
$ws.Range("D19").Value = "11.72"
$ws.Range("E19").Value = "  +1.31%  "
$ws.Range("D20").Value = "7.79"
$ws.Range("E20").Value = "  -2.95%  "
$ws.Range("D21").Value = "360.90"
$ws.Range("E21").Value = "  -2.30%  "
$ws.Range("D22").Value = "4.12"
$ws.Range("E22").Value = "  -2.55%  "
$ws.Range("E23").Value = "  -6.30%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").Value = "70.74"
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("D26").Value = "1.84"
$ws.Range("E26").Value = "  -5.44%  "
$ws.Range("E27").Value = "  -8.39%  "
$ws.Range("D28").Value = "0.995"
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("E29").Value = "  -2.32%  "
$ws.Range("E30").Value = "  -6.09%  "
$ws.Range("E31").Value = "  -2.69%  "
$ws.Range("D32").Value = "504.39"
$ws.Range("E32").Value = "  -6.92%  "
$ws.Range("E33").Value = "  -2.47%  "
$ws.Range("D34").Value = "1.27"
$ws.Range("E34").Value = "  -5.62%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").Value = "  -2.26%  "
$ws.Range("D37").Value = "159.28"
$ws.Range("E37").Value = "  +1.13%  "
$ws.Range("E38").Value = "  +1.53%  "
$ws.Range("E39").Value = "  -3.44%  "
$ws.Range("E40").Value = "  -0.83%  "
$ws.Range("E41").Value = "  -4.70%  "
$ws.Range("E42").Value = "  -5.13%  "
$ws.Range("E43").Value = "  -5.86%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "2.49"
$ws.Range("E45").Value = "  -3.13%  "
$ws.Range("E46").Value = "  -1.83%  "
$ws.Range("D47").Value = "141.93"
$ws.Range("D48").Value = "0.539"
$ws.Range("E48").Value = "  -4.55%  "
$ws.Range("E49").Value = "  -3.62%  "
$ws.Range("D50").Value = "0.0₆0266"
$ws.Range("E50").Value = "  -5.39%  "
$ws.Range("E51").Value = "  -4.19%  "

# Restore default (Normal) style on column D so no stray number-format
# style lingers on cells that did not need one.
$ws.Range("D2:D51").Style = "Normal"

